$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "PI hours" sheet: add a new "cfop" column (G) with per-person
#    cfop grouping, mirroring the existing "dept"/"app" columns.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("PI hours")

# Copy the header cell formatting (bold / centered / bordered) from
# the existing "app" header (F1) onto the new "cfop" header (G1).
$ws1.Range("F1").Copy()
$ws1.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws1.Range("G1").Value = "cfop"
$ws1.Range("G2").Value = "['cfop_MITRA']"
$ws1.Range("G3").Value = "['cfop_NH']"

# ------------------------------------------------------------------
# 2. Add a brand-new "cfop hours" sheet, after the existing sheets,
#    summarizing hours/percentage per cfop value.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "cfop hours"

# Match the look of the other summary sheets ("department hours",
# "unit(accumulative) hours"): bold/centered/bordered header row and
# a styled numeric index column.
$ws1.Range("B1:D1").Copy()
$ws4.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A2:A3").Copy()
$ws4.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws4.Range("B1").Value = "cfop"
$ws4.Range("C1").Value = "hours"
$ws4.Range("D1").Value = "percentage"

$ws4.Range("A2").Value = 0
$ws4.Range("B2").Value = "cfop_MITRA"
$ws4.Range("C2").Value = 14
$ws4.Range("D2").Value = 93.33333333333333

$ws4.Range("A3").Value = 1
$ws4.Range("B3").Value = "cfop_NH"
$ws4.Range("C3").Value = 1
$ws4.Range("D3").Value = 6.666666666666667

# Keep "PI hours" as the active/selected sheet, as it was before.
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
